$d = $word.ActiveDocument

# 1. Title heading
$d.Content.Find.Execute("Play Candy Bars Free - Review of IGT's Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Candy Bars Free - Colorful Slot Game with Wild Multipliers", 2)

# 2. "What we like" bullet list
$d.Content.Find.Execute("Colorful graphics style", $true, $false, $false, $false, $false, $true, 1, $false, "Colorful graphics and classic slot feel", 2)
$d.Content.Find.Execute("Classic slot feel", $true, $false, $false, $false, $false, $true, 1, $false, "Gumball Wilds with double and quadruple payouts", 2)
$d.Content.Find.Execute("Gumball Wilds double wins", $true, $false, $false, $false, $false, $true, 1, $false, "Wild Multipliers for increased winnings", 2)
$d.Content.Find.Execute("Frequent winning opportunities", $true, $false, $false, $false, $false, $true, 1, $false, "Progressive jackpots for additional excitement", 2)

# 3. "What we don't like" bullet list
$d.Content.Find.Execute("Limited number of symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Limited to 4 reels and 50 paylines. Less variety compared to other slots", 2)
$d.Content.Find.Execute("Missing more interactive features", $true, $false, $false, $false, $false, $true, 1, $false, "Lower maximum payout compared to some other slot games", 2)

# 4. Second occurrence of title text (bold text near end)
$d.Content.Find.Execute("Play Candy Bars Free - Review of IGT's Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Candy Bars Free - Colorful Slot Game with Wild Multipliers", 2)

# 5. Meta description (italic text)
$d.Content.Find.Execute("Read our review of Candy Bars by IGT. Play this colorful and classic slot game for free. Learn how to win Blackout Wins and Progressive Jackpots.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Candy Bars and play this colorful slot game for free. Enjoy Wild Multipliers and progressive jackpots.", 2)
